$d = $word.ActiveDocument

$replacements = @(
    @{old="60×48=2880"; new="37×34=1258"},
    @{old="35×52=1820"; new="92×78=7176"},
    @{old="35×24=840"; new="42×44=1848"},
    @{old="58×93=5394"; new="19×21=399"},
    @{old="60×45=2700"; new="29×88=2552"},
    @{old="48×16=768"; new="14×53=742"},
    @{old="97×78=7566"; new="13×46=598"},
    @{old="53×62=3286"; new="97×31=3007"},
    @{old="55×76=4180"; new="50×11=550"},
    @{old="96×61=5856"; new="33×95=3135"},
    @{old="79×55=4345"; new="55×81=4455"},
    @{old="89×12=1068"; new="85×98=8330"},
    @{old="40×36=1440"; new="79×62=4898"},
    @{old="36×36=1296"; new="21×17=357"},
    @{old="75×21=1575"; new="95×29=2755"},
    @{old="36×35=1260"; new="99×59=5841"},
    @{old="72×17=1224"; new="19×46=874"},
    @{old="75×91=6825"; new="27×94=2538"},
    @{old="23×88=2024"; new="22×77=1694"},
    @{old="55×63=3465"; new="40×20=800"},
    @{old="11×13=143"; new="97×44=4268"},
    @{old="17×54=918"; new="97×42=4074"},
    @{old="83×52=4316"; new="85×30=2550"},
    @{old="24×90=2160"; new="69×85=5865"},
    @{old="68×90=6120"; new="76×50=3800"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
